$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column J (new revenue component) and column AG (total) for rows 2-6
$ws.Range("J2").Value = 9093.950000000001
$ws.Range("AG2").Value = 75591.73

$ws.Range("J3").Value = 5820.5
$ws.Range("AG3").Value = 34160.05

$ws.Range("J4").Value = 1575.5
$ws.Range("AG4").Value = 28619.4

$ws.Range("J5").Value = 1824.02
$ws.Range("AG5").Value = 27365.05

$ws.Range("J6").Value = 18313.97
$ws.Range("AG6").Value = 165736.23
